# Applies the GSC export update:
# - Adds three new date rows (2025-11-18, 2025-11-19, 2025-11-20) to the "Chart" sheet
# - Updates the "Table" sheet's Validation status from "Started" to "Failed"

$wb = $excel.ActiveWorkbook

$chart = $wb.Worksheets.Item("Chart")
$table = $wb.Worksheets.Item("Table")

# Add the three new rows to the Chart sheet (rows 46-48)
# Force column A to be plain text so the date strings are not
# auto-converted into Excel date serial numbers.
$chart.Range("A46:A48").NumberFormat = "@"

$chart.Range("A46").Value = "2025-11-18"
$chart.Range("B46").Value = 23.0
$chart.Range("C46").Value = 1.0
$chart.Range("D46").Value = 0.0

$chart.Range("A47").Value = "2025-11-19"
$chart.Range("B47").Value = 23.0
$chart.Range("C47").Value = 1.0
$chart.Range("D47").Value = 0.0

$chart.Range("A48").Value = "2025-11-20"
$chart.Range("B48").Value = 23.0
$chart.Range("C48").Value = 1.0
$chart.Range("D48").Value = ""

# Update the Validation status in the Table sheet
$table.Range("B2").Value = "Failed"
